$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates, mirroring the scraped data refresh.

$ws.Range("D2").Value = "26.986.75"

$ws.Range("D3").Value = "1.649.58"
$ws.Range("E3").Value = "  +3.47%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.90"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +1.71%  "

$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  +3.76%  "

$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").Value = "1.882.85"
$ws.Range("E12").Value = "  +3.51%  "

$ws.Range("D13").Value = "1.657.18"

$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.14"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "26.967.01"
$ws.Range("E17").Value = "  +2.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.57"
$ws.Range("E18").Value = "  +3.71%  "

$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("E22").Value = "  +4.45%  "

$ws.Range("E23").Value = "  +4.44%  "

$ws.Range("E24").Value = "  +3.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  +1.94%  "

$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("E29").Value = "  +2.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("E32").Value = "  +3.19%  "

$ws.Range("D33").Value = "1.510.49"
$ws.Range("E33").Value = "  +2.67%  "

$ws.Range("E34").Value = "  +5.06%  "

$ws.Range("E35").Value = "  +9.10%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.574"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("E38").Value = "  +8.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.94"
$ws.Range("E40").Value = "  +3.43%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.76"
$ws.Range("E43").Value = "  +9.32%  "

$ws.Range("D44").Value = "1.789.97"
$ws.Range("E44").Value = "  +3.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("E45").Value = "  +2.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.917"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.35"
$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("E49").Value = "  +3.17%  "

$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("E51").Value = "  +2.07%  "
